$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "name" column header
$ws.Range("D1").Value = "name"

# Add name value for existing row 2
$ws.Range("D2").Value = "test_1"

# Duplicate row 2's test/data paths into a new row 3 for a second test config
$ws.Range("A3").Value = "test1/test_sheet.xlsx"
$ws.Range("B3").Value = "test1/data_sheet.xlsx"
$ws.Range("C3").Value = "http://localhost:9001"
$ws.Range("D3").Value = "Test_2"

# Hyperlink the new url cell just like C2
$ws.Hyperlinks.Add($ws.Range("C3"), "http://localhost:9001", "", "", "http://localhost:9001")

# Re-apply C2's cell formatting onto C3 (Hyperlinks.Add swaps in Excel's
# built-in Hyperlink style, but this column should keep the sheet's style)
$ws.Range("C2").Copy($ws.Range("C3"))

$excel.CutCopyMode = $false

$ws.Range("D3").Select()
